$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    'D2' = '37.489.77'
    'E2' = '  +5.28%  '
    'D3' = '2.058.30'
    'E3' = '  +3.72%  '
    'D5' = '253.56'
    'E5' = '  +2.95%  '
    'D6' = '0.655'
    'E6' = '  +3.01%  '
    'D7' = '67.71'
    'E7' = '  +15.41%  '
    'E8' = '  +0.00%  '
    'D9' = '0.386'
    'E9' = '  +6.64%  '
    'D10' = '59.82'
    'E10' = '  +1.89%  '
    'D11' = '0.0770'
    'E11' = '  +4.48%  '
    'E12' = '  +0.82%  '
    'D13' = '0.933'
    'E13' = '  -2.47%  '
    'D14' = '14.91'
    'E14' = '  +1.83%  '
    'D15' = '23.21'
    'E15' = '  +26.37%  '
    'D16' = '2.358.33'
    'D17' = '5.60'
    'E17' = '  +5.12%  '
    'D18' = '2.052.61'
    'E18' = '  +3.63%  '
    'D19' = '37.429.09'
    'E19' = '  +5.27%  '
    'D20' = '73.71'
    'E20' = '  +3.16%  '
    'D21' = '0.0₃0879'
    'E21' = '  +3.54%  '
    'D22' = '5.49'
    'D23' = '240.36'
    'E23' = '  +3.20%  '
    'E24' = '  +4.81%  '
    'E25' = '  +0.06%  '
    'D26' = '2.46'
    'E26' = '  +8.04%  '
    'D27' = '9.99'
    'E27' = '  +9.08%  '
    'D28' = '162.21'
    'E28' = '  -1.68%  '
    'D29' = '20.07'
    'E29' = '  +4.24%  '
    'E30' = '  +36.40%  '
    'E31' = '  +3.10%  '
    'D32' = '5.26'
    'E32' = '  +7.54%  '
    'E33' = '  +8.74%  '
    'D34' = '4.71'
    'E34' = '  +7.90%  '
    'D35' = '0.0630'
    'E35' = '  +5.74%  '
    'D36' = '2.45'
    'E36' = '  +0.68%  '
    'D37' = '6.16'
    'E37' = '  +15.20%  '
    'E38' = '  +0.07%  '
    'E39' = '  +4.01%  '
    'D40' = '3.17'
    'E40' = '  +39.37%  '
    'D41' = '0.103'
    'E41' = '  +14.72%  '
    'D42' = '1.27'
    'E42' = '  +3.53%  '
    'E43' = '  +5.79%  '
    'D44' = '17.52'
    'E44' = '  +8.24%  '
    'E45' = '  +5.40%  '
    'D46' = '0.0220'
    'E46' = '  +3.13%  '
    'D47' = '97.37'
    'E47' = '  +4.12%  '
    'D48' = '7.96'
    'E48' = '  +1.38%  '
    'D49' = '1.417.11'
    'E49' = '  +2.76%  '
    'E50' = '  +1.81%  '
    'D51' = '3.80'
    'E51' = '  +10.34%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    $cell.Style = 'Normal'
}
